$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 90, shifting rows 90:164 down to 91:165
$ws.Rows.Item(90).Insert()

# Fill the new row 90 with the new weekly data point
$ws.Cells.Item(90, 1).Value = 7
$ws.Cells.Item(90, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(90, 3).Value = "Ñuble"
$ws.Cells.Item(90, 4).Value = 44447
$ws.Cells.Item(90, 4).NumberFormat = $ws.Cells.Item(91, 4).NumberFormat
$ws.Cells.Item(90, 5).Value = 16
$ws.Cells.Item(90, 6).Value = 100114013
$ws.Cells.Item(90, 7).Value = "Zanahoria"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 300
$ws.Cells.Item(90, 11).Value = 5000
$ws.Cells.Item(90, 12).Value = 5500
$ws.Cells.Item(90, 13).Value = 5250
$ws.Cells.Item(90, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(90, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(90, 16).Value = 262
$ws.Cells.Item(90, 17).Value = 20
$ws.Cells.Item(90, 18).Value = "Hortaliza"
